$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Cells.Item(74, 8).Value = 4200.0586
$ws.Cells.Item(74, 9).Value = 4333.4165
$ws.Cells.Item(74, 10).Value = 3880
$ws.Cells.Item(74, 11).Value = 4333.4165
$ws.Cells.Item(74, 12).Value = 3880
$ws.Cells.Item(74, 13).Value = -3397.4165
$ws.Cells.Item(74, 14).Value = -5752

# Row 77
$ws.Cells.Item(77, 8).Value = 4200.0586
$ws.Cells.Item(77, 9).Value = 4333.4165
$ws.Cells.Item(77, 10).Value = 3880
$ws.Cells.Item(77, 11).Value = 21667.0825
$ws.Cells.Item(77, 12).Value = 19400
$ws.Cells.Item(77, 13).Value = -16987.0825
$ws.Cells.Item(77, 14).Value = -28760

# Row 129
$ws.Cells.Item(129, 8).Value = 439.16666
$ws.Cells.Item(129, 10).Value = 1200
$ws.Cells.Item(129, 12).Value = 3600
$ws.Cells.Item(129, 14).Value = -13600

# Row 132
$ws.Cells.Item(132, 8).Value = 743919.1
$ws.Cells.Item(132, 9).Value = 1338.7542
$ws.Cells.Item(132, 10).Value = 9803400
$ws.Cells.Item(132, 11).Value = 4016.2626
$ws.Cells.Item(132, 12).Value = 29410200
$ws.Cells.Item(132, 13).Value = -1486.2626
$ws.Cells.Item(132, 14).Value = -29415260

$ws = $wb.Worksheets.Item("ARM")
# Row 39
$ws.Cells.Item(39, 8).Value = 4804.25
$ws.Cells.Item(39, 9).Value = 3072
$ws.Cells.Item(39, 10).Value = 10001
$ws.Cells.Item(39, 11).Value = 3072
$ws.Cells.Item(39, 12).Value = 10001
$ws.Cells.Item(39, 13).Value = -2552
$ws.Cells.Item(39, 14).Value = -11041

# Row 74
$ws.Cells.Item(74, 8).Value = 6632924
$ws.Cells.Item(74, 9).Value = 8367790.5
$ws.Cells.Item(74, 10).Value = 127175
$ws.Cells.Item(74, 11).Value = 8367790.5
$ws.Cells.Item(74, 12).Value = 127175
$ws.Cells.Item(74, 13).Value = -8366916.5
$ws.Cells.Item(74, 14).Value = -128923

# Row 77
$ws.Cells.Item(77, 8).Value = 6632924
$ws.Cells.Item(77, 9).Value = 8367790.5
$ws.Cells.Item(77, 10).Value = 127175
$ws.Cells.Item(77, 11).Value = 41838952.5
$ws.Cells.Item(77, 12).Value = 635875
$ws.Cells.Item(77, 13).Value = -41834584.5
$ws.Cells.Item(77, 14).Value = -644611

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Cells.Item(62, 8).Value = 30995
$ws.Cells.Item(62, 10).Value = 30995
$ws.Cells.Item(62, 12).Value = 30995
$ws.Cells.Item(62, 14).Value = -32367

# Row 65
$ws.Cells.Item(65, 8).Value = 30995
$ws.Cells.Item(65, 10).Value = 30995
$ws.Cells.Item(65, 12).Value = 92985
$ws.Cells.Item(65, 14).Value = -99849

# Row 86
$ws.Cells.Item(86, 8).Value = 10519.63
$ws.Cells.Item(86, 9).Value = 11783.305
$ws.Cells.Item(86, 10).Value = 3253.5
$ws.Cells.Item(86, 11).Value = 11783.305
$ws.Cells.Item(86, 12).Value = 3253.5
$ws.Cells.Item(86, 13).Value = -10660.305
$ws.Cells.Item(86, 14).Value = -5499.5

# Row 89
$ws.Cells.Item(89, 8).Value = 10519.63
$ws.Cells.Item(89, 9).Value = 11783.305
$ws.Cells.Item(89, 10).Value = 3253.5
$ws.Cells.Item(89, 11).Value = 58916.525
$ws.Cells.Item(89, 12).Value = 16267.5
$ws.Cells.Item(89, 13).Value = -53300.525
$ws.Cells.Item(89, 14).Value = -27499.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1320
$ws.Cells.Item(16, 9).Value = 1024.5
$ws.Cells.Item(16, 10).Value = 1779.6666
$ws.Cells.Item(16, 11).Value = 1024.5
$ws.Cells.Item(16, 12).Value = 1779.6666
$ws.Cells.Item(16, 13).Value = -737.5
$ws.Cells.Item(16, 14).Value = -2353.6666

# Row 35
$ws.Cells.Item(35, 8).Value = 800
$ws.Cells.Item(35, 9).Value = 800
$ws.Cells.Item(35, 11).Value = 800
$ws.Cells.Item(35, 13).Value = -506

# Row 58
$ws.Cells.Item(58, 8).Value = 17242668
$ws.Cells.Item(58, 9).Value = 19609022
$ws.Cells.Item(58, 10).Value = 2081.5715
$ws.Cells.Item(58, 11).Value = 19609022
$ws.Cells.Item(58, 12).Value = 2081.5715
$ws.Cells.Item(58, 13).Value = -19608819
$ws.Cells.Item(58, 14).Value = -2487.5715

# Row 105
$ws.Cells.Item(105, 8).Value = 1847.9333
$ws.Cells.Item(105, 9).Value = 1809
$ws.Cells.Item(105, 10).Value = 2003.6666
$ws.Cells.Item(105, 11).Value = 1809
$ws.Cells.Item(105, 12).Value = 2003.6666
$ws.Cells.Item(105, 13).Value = -62
$ws.Cells.Item(105, 14).Value = -5497.6666

# Row 113
$ws.Cells.Item(113, 8).Value = 1320
$ws.Cells.Item(113, 9).Value = 1024.5
$ws.Cells.Item(113, 10).Value = 1779.6666
$ws.Cells.Item(113, 11).Value = 1024.5
$ws.Cells.Item(113, 12).Value = 1779.6666
$ws.Cells.Item(113, 13).Value = 1145.5
$ws.Cells.Item(113, 14).Value = -6119.6666

# Row 132
$ws.Cells.Item(132, 8).Value = 44026.637
$ws.Cells.Item(132, 9).Value = 25787.121
$ws.Cells.Item(132, 10).Value = 168663.33
$ws.Cells.Item(132, 11).Value = 77361.363
$ws.Cells.Item(132, 12).Value = 505989.99
$ws.Cells.Item(132, 13).Value = -74831.363
$ws.Cells.Item(132, 14).Value = -511049.99

# Row 134
$ws.Cells.Item(134, 8).Value = 56057.7
$ws.Cells.Item(134, 9).Value = 2249.2307
$ws.Cells.Item(134, 10).Value = 155987.72
$ws.Cells.Item(134, 11).Value = 6747.6921
$ws.Cells.Item(134, 12).Value = 467963.16
$ws.Cells.Item(134, 13).Value = -4212.6921
$ws.Cells.Item(134, 14).Value = -473033.16

# Row 136
$ws.Cells.Item(136, 8).Value = 17242668
$ws.Cells.Item(136, 9).Value = 19609022
$ws.Cells.Item(136, 10).Value = 2081.5715
$ws.Cells.Item(136, 11).Value = 58827066
$ws.Cells.Item(136, 12).Value = 6244.7145
$ws.Cells.Item(136, 13).Value = -58824516
$ws.Cells.Item(136, 14).Value = -11344.7145

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3080.7073
$ws.Cells.Item(80, 9).Value = 2500.6
$ws.Cells.Item(80, 10).Value = 3415.3845
$ws.Cells.Item(80, 11).Value = 2500.6
$ws.Cells.Item(80, 12).Value = 3415.3845
$ws.Cells.Item(80, 13).Value = -1502.6
$ws.Cells.Item(80, 14).Value = -5411.3845

# Row 83
$ws.Cells.Item(83, 8).Value = 3080.7073
$ws.Cells.Item(83, 9).Value = 2500.6
$ws.Cells.Item(83, 10).Value = 3415.3845
$ws.Cells.Item(83, 11).Value = 12503
$ws.Cells.Item(83, 12).Value = 17076.9225
$ws.Cells.Item(83, 13).Value = -7511
$ws.Cells.Item(83, 14).Value = -27060.9225

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Cells.Item(93, 8).Value = 1393
$ws.Cells.Item(93, 9).Value = 1393
$ws.Cells.Item(93, 11).Value = 1393
$ws.Cells.Item(93, 13).Value = -145

# Row 100
$ws.Cells.Item(100, 8).Value = 1340.7778
$ws.Cells.Item(100, 9).Value = 1152.4286
$ws.Cells.Item(100, 11).Value = 1152.4286
$ws.Cells.Item(100, 13).Value = -611.4286

# Row 132
$ws.Cells.Item(132, 8).Value = 271726.22
$ws.Cells.Item(132, 9).Value = 314000.12
$ws.Cells.Item(132, 10).Value = 204088
$ws.Cells.Item(132, 11).Value = 942000.36
$ws.Cells.Item(132, 12).Value = 612264
$ws.Cells.Item(132, 13).Value = -939470.36
$ws.Cells.Item(132, 14).Value = -617324

# Row 136
$ws.Cells.Item(136, 8).Value = 70258
$ws.Cells.Item(136, 9).Value = 62000.61
$ws.Cells.Item(136, 11).Value = 186001.83
$ws.Cells.Item(136, 13).Value = -183451.83

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 66787.32000000001
$ws.Cells.Item(132, 9).Value = 50800.35
$ws.Cells.Item(132, 10).Value = 158141.42
$ws.Cells.Item(132, 11).Value = 152401.05
$ws.Cells.Item(132, 12).Value = 474424.26
$ws.Cells.Item(132, 13).Value = -149871.05
$ws.Cells.Item(132, 14).Value = -479484.26

# Row 136
$ws.Cells.Item(136, 8).Value = 35484.535
$ws.Cells.Item(136, 9).Value = 20900.06
$ws.Cells.Item(136, 10).Value = 114888.89
$ws.Cells.Item(136, 11).Value = 62700.18000000001
$ws.Cells.Item(136, 12).Value = 344666.67
$ws.Cells.Item(136, 13).Value = -60150.18000000001
$ws.Cells.Item(136, 14).Value = -349766.67
